$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(60, 1).Value = "2025/12/05 05:00"
$ws.Cells.Item(60, 2).Value = "-"
$ws.Cells.Item(60, 3).Value = "-"
$ws.Cells.Item(60, 4).Value = "-"
$ws.Cells.Item(60, 5).Value = "-"
$ws.Cells.Item(60, 6).Value = "-"
$ws.Cells.Item(60, 7).Value = "-"
